$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: updated student IDs ---
$ws.Range("A2").Value = 2020211033376
$ws.Range("A3").Value = 2020211033378
$ws.Range("A4").Value = 2020211023021
$ws.Range("A5").Value = 2020211023022
$ws.Range("A6").Value = 2020211023023

# --- Column B: updated student names ---
# (order matters for how new strings land in the shared-strings table,
#  so write B2, B3, B4, B6, B5 to match the saved workbook's string order)
$ws.Range("B2").Value = "zhangsan"
$ws.Range("B3").Value = "深深的"
$ws.Range("B4").Value = "得到的"
$ws.Range("B6").Value = "猜猜猜"
$ws.Range("B5").Value = "从吃"

# --- Column H: all rows now reference the same supervising teacher ---
$ws.Range("H2").Value = "王老师"
$ws.Range("H3").Value = "王老师"
$ws.Range("H4").Value = "王老师"
$ws.Range("H5").Value = "王老师"
$ws.Range("H6").Value = "王老师"

# --- Update the saved selection to match the latest user action ---
[void]$ws.Range("H2:H6").Select()
